$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Bring the "Repayment schedule" sheet to the front (it becomes the active tab,
# and the previously active "Transactions" sheet loses the active-tab flag).
$ws.Activate()

# Insert a new (blank-header) column before the old "Late" column (N), shifting
# "Late", "heading" and "Outstanding" one column to the right.
$ws.Columns("N:N").Insert()

# Resize the "In Advance" column and the freshly inserted blank column together
# to match the width used by the "Outstanding" column.
$ws.Columns("M:N").ColumnWidth = 11.85546875

# Restore the selection on the sheet to the cell the author ended up on.
$ws.Range("S8").Select()
